$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44162
$ws.Range("K2").Value = 'Castle Brite'
$ws.Range("M2").Value = 70
$ws.Range("N2").Value = 8500
$ws.Range("O2").Value = 8500
$ws.Range("P2").Value = 8500
$ws.Range("Q2").Value = '$/bandeja 10 kilos'
$ws.Range("S2").Value = 850
$ws.Range("T2").Value = 10

# Row 3
$ws.Range("D3").Value = 44162
$ws.Range("K3").Value = 'Castle Brite'
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 75
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14400
$ws.Range("Q3").Value = '$/caja 18 kilos'
$ws.Range("S3").Value = 800
$ws.Range("T3").Value = 18

# Row 4
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("S4").Value = 1000

# Row 5
$ws.Range("D5").Value = 44187
$ws.Range("K5").Value = 'Patterson'
$ws.Range("L5").Value = 'Segunda'
$ws.Range("M5").Value = 95
$ws.Range("N5").Value = 13500
$ws.Range("O5").Value = 13500
$ws.Range("P5").Value = 13500
$ws.Range("Q5").Value = '$/caja 15 kilos granel'
$ws.Range("S5").Value = 900
$ws.Range("T5").Value = 15

# Row 6
$ws.Range("D6").Value = 44187
$ws.Range("K6").Value = 'Patterson'
$ws.Range("L6").Value = 'Tercera'
$ws.Range("M6").Value = 120
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 12000
$ws.Range("Q6").Value = '$/caja 15 kilos granel'
$ws.Range("S6").Value = 800
$ws.Range("T6").Value = 15

# Row 7
$ws.Range("D7").Value = 44523
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 320
$ws.Range("N7").Value = 10000
$ws.Range("O7").Value = 10000
$ws.Range("P7").Value = 10000
$ws.Range("Q7").Value = '$/bandeja 10 kilos'
$ws.Range("S7").Value = 1000
$ws.Range("T7").Value = 10

# Row 8
$ws.Range("D8").Value = 44175
$ws.Range("K8").Value = 'Modesto'
$ws.Range("M8").Value = 140
$ws.Range("N8").Value = 11000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 11571
$ws.Range("Q8").Value = '$/caja 12 kilos'
$ws.Range("S8").Value = 964
$ws.Range("T8").Value = 12

# Row 9
$ws.Range("D9").Value = 44169
$ws.Range("K9").Value = 'Dina'
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 10000
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 1000

# Row 10
$ws.Range("D10").Value = 44194
$ws.Range("K10").Value = 'Patterson'
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 120
$ws.Range("N10").Value = 13000
$ws.Range("O10").Value = 13000
$ws.Range("P10").Value = 13000
$ws.Range("Q10").Value = '$/caja 15 kilos'
$ws.Range("S10").Value = 867
$ws.Range("T10").Value = 15

# Row 11
$ws.Range("D11").Value = 44174
$ws.Range("K11").Value = 'Modesto'
$ws.Range("M11").Value = 120
$ws.Range("N11").Value = 8500
$ws.Range("O11").Value = 8500
$ws.Range("P11").Value = 8500
$ws.Range("Q11").Value = '$/bandeja 10 kilos'
$ws.Range("R11").Value = 'Región Metropolitana'
$ws.Range("S11").Value = 850
$ws.Range("T11").Value = 10

# Row 12
$ws.Range("D12").Value = 44174
$ws.Range("M12").Value = 180
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 15000
$ws.Range("Q12").Value = '$/caja 18 kilos'
$ws.Range("R12").Value = 'Región Metropolitana'
$ws.Range("S12").Value = 833
$ws.Range("T12").Value = 18

# Row 13
$ws.Range("D13").Value = 44174
$ws.Range("K13").Value = 'Modesto'
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 120
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 12000
$ws.Range("Q13").Value = '$/caja 18 kilos'
$ws.Range("R13").Value = 'Región Metropolitana'
$ws.Range("S13").Value = 667
$ws.Range("T13").Value = 18

# Row 14
$ws.Range("D14").Value = 44195
$ws.Range("K14").Value = 'Patterson'
$ws.Range("M14").Value = 124
$ws.Range("N14").Value = 13000
$ws.Range("O14").Value = 13000
$ws.Range("P14").Value = 13000
$ws.Range("Q14").Value = '$/caja 15 kilos'
$ws.Range("S14").Value = 867
$ws.Range("T14").Value = 15

# Row 15
$ws.Range("D15").Value = 44176
$ws.Range("K15").Value = 'Modesto'
$ws.Range("M15").Value = 115
$ws.Range("N15").Value = 11000
$ws.Range("O15").Value = 12000
$ws.Range("P15").Value = 11609
$ws.Range("Q15").Value = '$/caja 12 kilos'
$ws.Range("S15").Value = 967
$ws.Range("T15").Value = 12

# Row 16
$ws.Range("D16").Value = 44159
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 80
$ws.Range("N16").Value = 8000
$ws.Range("O16").Value = 8000
$ws.Range("P16").Value = 8000
$ws.Range("S16").Value = 800

# Row 17
$ws.Range("D17").Value = 44159
$ws.Range("L17").Value = 'Segunda'
$ws.Range("M17").Value = 65
$ws.Range("N17").Value = 7000
$ws.Range("O17").Value = 7000
$ws.Range("P17").Value = 7000
$ws.Range("S17").Value = 700

# Row 18
$ws.Range("D18").Value = 44519
$ws.Range("M18").Value = 300
$ws.Range("N18").Value = 22500
$ws.Range("O18").Value = 22500
$ws.Range("P18").Value = 22500
$ws.Range("Q18").Value = '$/caja 15 kilos granel'
$ws.Range("R18").Value = 'Región de O''Higgins'
$ws.Range("S18").Value = 1500
$ws.Range("T18").Value = 15

# Row 19
$ws.Range("D19").Value = 44168
$ws.Range("K19").Value = 'Dina'
$ws.Range("L19").Value = 'Especial'
$ws.Range("M19").Value = 40
$ws.Range("N19").Value = 14000
$ws.Range("O19").Value = 14000
$ws.Range("P19").Value = 14000
$ws.Range("R19").Value = 'Región de O''Higgins'
$ws.Range("S19").Value = 1400

# Row 20
$ws.Range("D20").Value = 44166
$ws.Range("K20").Value = 'Castle Brite'
$ws.Range("M20").Value = 120
$ws.Range("N20").Value = 10000
$ws.Range("O20").Value = 10000
$ws.Range("P20").Value = 10000
$ws.Range("Q20").Value = '$/bandeja 10 kilos'
$ws.Range("S20").Value = 1000
$ws.Range("T20").Value = 10

# Row 21
$ws.Range("D21").Value = 44166
$ws.Range("K21").Value = 'Castle Brite'
$ws.Range("N21").Value = 8000
$ws.Range("O21").Value = 8000
$ws.Range("P21").Value = 8000
$ws.Range("Q21").Value = '$/bandeja 10 kilos'
$ws.Range("S21").Value = 800
$ws.Range("T21").Value = 10

# Row 22
$ws.Range("D22").Value = 44525
$ws.Range("K22").Value = 'Castle Brite'
$ws.Range("L22").Value = 'Especial'
$ws.Range("M22").Value = 300
$ws.Range("N22").Value = 25200
$ws.Range("O22").Value = 25200
$ws.Range("P22").Value = 25200
$ws.Range("S22").Value = 1400

# Row 23
$ws.Range("D23").Value = 44525
$ws.Range("M23").Value = 250
$ws.Range("N23").Value = 21600
$ws.Range("O23").Value = 21600
$ws.Range("P23").Value = 21600
$ws.Range("Q23").Value = '$/caja 18 kilos'
$ws.Range("S23").Value = 1200
$ws.Range("T23").Value = 18

# Row 24
$ws.Range("D24").Value = 44160
$ws.Range("K24").Value = 'Castle Brite'
$ws.Range("M24").Value = 25
$ws.Range("N24").Value = 8000
$ws.Range("O24").Value = 8000
$ws.Range("P24").Value = 8000
$ws.Range("R24").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S24").Value = 800

# Row 25
$ws.Range("D25").Value = 44160
$ws.Range("K25").Value = 'Castle Brite'
$ws.Range("L25").Value = 'Segunda'
$ws.Range("M25").Value = 40
$ws.Range("N25").Value = 7000
$ws.Range("O25").Value = 7000
$ws.Range("P25").Value = 7000
$ws.Range("Q25").Value = '$/bandeja 10 kilos'
$ws.Range("S25").Value = 700
$ws.Range("T25").Value = 10

# Row 26
$ws.Range("D26").Value = 44167
$ws.Range("K26").Value = 'Castle Brite'
$ws.Range("L26").Value = 'Especial'
$ws.Range("M26").Value = 85
$ws.Range("N26").Value = 10000
$ws.Range("O26").Value = 10000
$ws.Range("P26").Value = 10000
$ws.Range("Q26").Value = '$/bandeja 10 kilos'
$ws.Range("S26").Value = 1000
$ws.Range("T26").Value = 10

# Row 27
$ws.Range("D27").Value = 44167
$ws.Range("N27").Value = 9500
$ws.Range("O27").Value = 9500
$ws.Range("P27").Value = 9500
$ws.Range("R27").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S27").Value = 950

# Row 28
$ws.Range("D28").Value = 44167
$ws.Range("L28").Value = 'Primera'
$ws.Range("M28").Value = 60
$ws.Range("N28").Value = 15000
$ws.Range("O28").Value = 15000
$ws.Range("P28").Value = 15000
$ws.Range("Q28").Value = '$/caja 18 kilos'
$ws.Range("R28").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S28").Value = 833
$ws.Range("T28").Value = 18

# Row 29
$ws.Range("D29").Value = 44189
$ws.Range("M29").Value = 130
$ws.Range("N29").Value = 12000
$ws.Range("O29").Value = 12000
$ws.Range("P29").Value = 12000
$ws.Range("Q29").Value = '$/caja 18 kilos'
$ws.Range("S29").Value = 667
$ws.Range("T29").Value = 18
